$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Data" section paragraph - rewrite text (FRED / Macrotrends description)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "We pulled from two main sources for all of the data, the Federal Reserve Bank of St. Louis (FRED) and Macrotrends (website).  FRED has lots of valuable datasets for various economic and demographic statistics.  We used Macrotrends to get the oil price history table.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We pulled the data from two main sources, the Federal Reserve Bank of St. Louis (FRED) and Macrotrends (website).  FRED has numerous valuable datasets consisting of various economic and demographic statistics.  Macrotrends provides WTI oil price history.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Extraction paragraph - merge "The "+"Macrotrends"+" data " into a single
#    run and drop the spell-check proofErr wrapper around "Macrotrends".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The Macrotrends data (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The Macrotrends data (",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Transform paragraph - rewrite text
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The data transformation include aggregation (to 1 row per year), converting some #'s from thousands to millions (population) and renaming columns/indexes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The data transformation includes averaging monthly and quarterly values per year to derive an annual value of the statistics, converting all numbers into the true value as some datasets, population, reported the data as per thousand and renaming columns/indexes.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Flask API paragraph - rewrite text (keep the trailing two-space run as-is)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Our Flask program, app.py, can be used to pull a combined dataset (all five tables), or each table separately.  The six calls can be accessed from the local host address on any browser.  The index.html site lists out the addresses for the calls.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Our Flask program, app.py, provides a combined dataset (all five tables), or each table separately.  The six calls can be accessed from the local host address on any browser.  The index.html site lists out the addresses for the calls for the api calls.",
    2) | Out-Null

Write-Output "text rewrites done"

# ---------------------------------------------------------------------------
# 5. Insert new "Normalization of Datasets" heading + paragraph right after
#    the last bullet of "Potential Uses" ("Has the impact ... since 2000?").
#    Also relocate the "_GoBack" bookmark into the new paragraph.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete() | Out-Null

$anchorPar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Has the impact of oil price on employment levels*") {
        $anchorPar = $d.Paragraphs($i)
        break
    }
}

$anchorPar.Range.InsertParagraphAfter() | Out-Null
$headingPar = $anchorPar.Next()
$headingPar.Range.ListFormat.RemoveNumbers()
$headingPar.Range.Style = "Heading 1"
$headingPar.Range.Text = "Normalization of Datasets"

$headingPar.Range.InsertParagraphAfter() | Out-Null
$bodyPar = $headingPar.Next()
$bodyPar.Range.ListFormat.RemoveNumbers()
$bodyPar.Range.Style = "Normal"

$insPt = $d.Range($bodyPar.Range.Start, $bodyPar.Range.Start)
$insPt.InsertAfter("The population and employment datasets could be normalized by dividing them by a thousand as they were initially reported.  We however elected not to do so as the datasets were small and doing so would require more documentation concerning the datasets available via the api.") | Out-Null

$spacePt = $d.Range($bodyPar.Range.End - 1, $bodyPar.Range.End - 1)
$spacePt.InsertAfter(" ") | Out-Null

$bmLoc = $bodyPar.Range.End - 2
$bmRange = $d.Range($bmLoc, $bmLoc)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$bodyPar.Range.InsertParagraphAfter() | Out-Null
$blankPar = $bodyPar.Next()
$blankPar.Range.ListFormat.RemoveNumbers()
$blankPar.Range.Style = "Normal"

Write-Output "normalization section inserted"
